$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Row 3: PERIOD TO EXPIRE -98 -> -99, LAST UPDATE 03-Nov-2025 -> 04-Nov-2025
$ws.Range("H3").Value = -99

$i3 = $ws.Range("I3")
$i3.Formula = '="04-Nov-2025"'
$i3.Copy()
$i3.PasteSpecial(-4163)

# Row 4: PERIOD TO EXPIRE 700 -> 699, LAST UPDATE 03-Nov-2025 -> 04-Nov-2025
$ws.Range("H4").Value = 699

$i4 = $ws.Range("I4")
$i4.Formula = '="04-Nov-2025"'
$i4.Copy()
$i4.PasteSpecial(-4163)

$excel.CutCopyMode = 0
